# Refresh the cryptos price/volume table (Price column D, Volume(1h) column E).
# Values that look like plain numbers (single decimal point, e.g. "581.98")
# are prefixed with a leading apostrophe so Excel keeps them as text (quote
# prefix) instead of silently converting them to numeric cells, matching the
# original inline-string formatting of the sheet. Values that already read
# as non-numeric text (e.g. "65.945.10", which has two dots) need no prefix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.945.10"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "3.472.41"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'581.98"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").Value = "'173.15"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("D9").Value = "3.470.34"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("E10").Value = "  -5.78%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("D12").Value = "'0.410"
$ws.Range("E12").Value = "  -3.89%  "
$ws.Range("D13").Value = "4.067.34"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "'29.79"
$ws.Range("E15").Value = "  -6.73%  "
$ws.Range("D16").Value = "66.042.73"
$ws.Range("E16").Value = "  -2.34%  "
$ws.Range("D17").Value = "'0.0000171"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").Value = "3.462.47"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("D20").Value = "'13.85"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").Value = "'366.28"
$ws.Range("E21").Value = "  -5.89%  "
$ws.Range("D22").Value = "'7.72"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "'72.11"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").Value = "'0.534"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").Value = "'0.0000124"
$ws.Range("E26").Value = "  +4.17%  "
$ws.Range("D27").Value = "'9.56"
$ws.Range("E27").Value = "  -6.46%  "
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").Value = "'23.84"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").Value = "'5.75"
$ws.Range("E31").Value = "  -4.95%  "
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -6.78%  "
$ws.Range("D35").Value = "'7.10"
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("D37").Value = "'159.26"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").Value = "'29.19"
$ws.Range("E38").Value = "  +13.15%  "
$ws.Range("D39").Value = "'0.890"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").Value = "2.798.81"
$ws.Range("E40").Value = "  +3.89%  "
$ws.Range("E41").Value = "  -5.38%  "
$ws.Range("E42").Value = "  -6.69%  "
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("D44").Value = "'4.43"
$ws.Range("E44").Value = "  -3.59%  "
$ws.Range("D45").Value = "'0.0681"
$ws.Range("E45").Value = "  -4.56%  "
$ws.Range("D46").Value = "'40.05"
$ws.Range("E46").Value = "  -2.65%  "
$ws.Range("D47").Value = "'24.13"
$ws.Range("E47").Value = "  -7.00%  "
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("D49").Value = "'312.88"
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("D50").Value = "'0.815"
$ws.Range("E50").Value = "  -2.42%  "
$ws.Range("D51").Value = "'0.101"
$ws.Range("E51").Value = "  -2.86%  "